$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D, E, G columns for data rows (2-51) are formatted as Text so
# values like "332.92" or "2" are stored as strings, matching the source data.
$ws.Range("D2:E51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

$ws.Range("D2").Value = "332.92"
$ws.Range("E2").Value = "2.21%"
$ws.Range("G2").Value = "2"

$ws.Range("D3").Value = "45.90"
$ws.Range("E3").Value = "4.55%"
$ws.Range("G3").Value = "2"

$ws.Range("D4").Value = "5.654"
$ws.Range("E4").Value = "2.89%"
$ws.Range("G4").Value = "2"

$ws.Range("D5").Value = "0.08370"
$ws.Range("E5").Value = "4.73%"
$ws.Range("G5").Value = "2"

$ws.Range("D6").Value = "2.050"
$ws.Range("E6").Value = "3.41%"
$ws.Range("G6").Value = "2"

$ws.Range("D7").Value = "0.9810"
$ws.Range("E7").Value = "3.82%"
$ws.Range("G7").Value = "2"

$ws.Range("D8").Value = "0.1144"
$ws.Range("E8").Value = "0.08%"
$ws.Range("G8").Value = "2"

$ws.Range("D9").Value = "0.1947"
$ws.Range("E9").Value = "5.80%"
$ws.Range("G9").Value = "2"

$ws.Range("E10").Value = "-12.86%"
$ws.Range("G10").Value = "2"

$ws.Range("D11").Value = "0.1003"
$ws.Range("E11").Value = "4.95%"
$ws.Range("G11").Value = "2"

$ws.Range("D12").Value = "0.04604"
$ws.Range("E12").Value = "-3.09%"
$ws.Range("G12").Value = "2"

$ws.Range("D13").Value = "0.1058"
$ws.Range("E13").Value = "-0.65%"
$ws.Range("G13").Value = "2"

$ws.Range("D14").Value = "0.001285"
$ws.Range("E14").Value = "1.65%"
$ws.Range("G14").Value = "2"

$ws.Range("D15").Value = "0.006049"
$ws.Range("E15").Value = "5.65%"
$ws.Range("G15").Value = "2"

$ws.Range("D16").Value = "3.376"
$ws.Range("E16").Value = "0.10%"
$ws.Range("G16").Value = "2"

$ws.Range("D17").Value = "4.458"
$ws.Range("E17").Value = "3.70%"
$ws.Range("G17").Value = "2"

$ws.Range("D18").Value = "2.612"
$ws.Range("E18").Value = "1.92%"
$ws.Range("G18").Value = "2"

$ws.Range("D19").Value = "0.3348"
$ws.Range("E19").Value = "-3.77%"
$ws.Range("G19").Value = "2"

$ws.Range("D20").Value = "0.1380"
$ws.Range("E20").Value = "-1.78%"
$ws.Range("G20").Value = "2"

$ws.Range("D21").Value = "0.2592"
$ws.Range("E21").Value = "1.79%"
$ws.Range("G21").Value = "2"

$ws.Range("D22").Value = "0.04109"
$ws.Range("E22").Value = "1.24%"
$ws.Range("G22").Value = "2"

$ws.Range("D23").Value = "0.001306"
$ws.Range("E23").Value = "4.95%"
$ws.Range("G23").Value = "2"

$ws.Range("D24").Value = "0.004425"
$ws.Range("E24").Value = "3.04%"
$ws.Range("G24").Value = "2"

$ws.Range("D25").Value = "0.0001280"
$ws.Range("E25").Value = "7.54%"
$ws.Range("G25").Value = "2"

$ws.Range("D26").Value = "0.0003740"
$ws.Range("E26").Value = "-0.12%"
$ws.Range("G26").Value = "2"

$ws.Range("G27").Value = "2"

$ws.Range("G28").Value = "2"

$ws.Range("G29").Value = "2"

$ws.Range("G30").Value = "2"

$ws.Range("G31").Value = "2"

$ws.Range("G32").Value = "2"

$ws.Range("G33").Value = "2"

$ws.Range("G34").Value = "2"

$ws.Range("G35").Value = "2"

$ws.Range("G36").Value = "2"

$ws.Range("G37").Value = "2"

$ws.Range("D38").Value = "0.02785"
$ws.Range("E38").Value = "10.72%"
$ws.Range("G38").Value = "2"

$ws.Range("D39").Value = "0.05804"
$ws.Range("E39").Value = "5.33%"
$ws.Range("G39").Value = "2"

$ws.Range("D40").Value = "0.007701"
$ws.Range("E40").Value = "2.24%"
$ws.Range("G40").Value = "2"

$ws.Range("D41").Value = "0.1439"
$ws.Range("E41").Value = "3.61%"
$ws.Range("G41").Value = "2"

$ws.Range("D42").Value = "0.007197"
$ws.Range("E42").Value = "-2.56%"
$ws.Range("G42").Value = "2"

$ws.Range("E43").Value = "-2.09%"
$ws.Range("G43").Value = "2"

$ws.Range("D44").Value = "0.008426"
$ws.Range("E44").Value = "0.53%"
$ws.Range("G44").Value = "2"

$ws.Range("D45").Value = "0.00007153"
$ws.Range("E45").Value = "0.71%"
$ws.Range("G45").Value = "2"

$ws.Range("D46").Value = "0.00000000750"
$ws.Range("E46").Value = "-0.01%"
$ws.Range("G46").Value = "2"

$ws.Range("D47").Value = "0.0005800"
$ws.Range("E47").Value = "-0.20%"
$ws.Range("G47").Value = "2"

$ws.Range("D48").Value = "0.003478"
$ws.Range("E48").Value = "-1.75%"
$ws.Range("G48").Value = "2"

$ws.Range("D49").Value = "0.003498"
$ws.Range("E49").Value = "51.96%"
$ws.Range("G49").Value = "2"

$ws.Range("D50").Value = "0.00002100"
$ws.Range("E50").Value = "-0.01%"
$ws.Range("G50").Value = "2"

$ws.Range("D51").Value = "0.0002000"
$ws.Range("E51").Value = "-0.01%"
$ws.Range("G51").Value = "2"
